# fix for tables query to .xlsx
#
# Target end-state (per diff):
#   Sheets in order: answers, hospitals, icu_beds_available, patient_responses,
#                    patients, questions
#   - "answers"  : unchanged header row (id, question_id, answers, weight)
#   - "hospitals": new sheet (id, hospital_name, total_icu_beds, total_vents,
#                  working_vents, location, latitutde, longitude)
#   - "icu_beds_available": new sheet (id, hospital_id, beds_in_use,
#                  vents_in_use, recorded_at) + one data row, recorded_at
#                  formatted as a date/time (numFmt "yyyy-mm-dd h:mm:ss")
#   - "patient_responses": new sheet (id, patient_id, question_id, answer_id,
#                  recorded_at)
#   - "patients" : new sheet (id, hospital_id, patient_name, patient_mr_no,
#                  entry_point, recorded_at)
#   - "questions": the old default "Sheet", renamed, with header row
#                  (id, question, question_type, weight)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the four brand-new sheets. Worksheets.Add() always inserts the
#    new sheet at the very front of the tab strip, so adding them in the
#    reverse of the desired final order leaves them front-to-back in the
#    right order once we're done.
# ---------------------------------------------------------------------------
$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item(1).Name = "patients"

$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item(1).Name = "patient_responses"

$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item(1).Name = "icu_beds_available"

$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item(1).Name = "hospitals"

# Current order: hospitals, icu_beds_available, patient_responses, patients, Sheet, answers

# ---------------------------------------------------------------------------
# 2. Put the pre-existing "answers" tab at the very front, and rename the
#    original default "Sheet" tab to "questions" (it is already last).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("answers").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("Sheet").Name = "questions"

# Final order should now be:
# answers, hospitals, icu_beds_available, patient_responses, patients, questions
$wb.Worksheets.Item("questions").Move($null, $wb.Worksheets.Item("patients"))

# ---------------------------------------------------------------------------
# 3. Populate "hospitals"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("hospitals")
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "hospital_name"
$ws.Cells.Item(1,3).Value = "total_icu_beds"
$ws.Cells.Item(1,4).Value = "total_vents"
$ws.Cells.Item(1,5).Value = "working_vents"
$ws.Cells.Item(1,6).Value = "location"
$ws.Cells.Item(1,7).Value = "latitutde"
$ws.Cells.Item(1,8).Value = "longitude"

# ---------------------------------------------------------------------------
# 4. Populate "icu_beds_available" (header + one data row; recorded_at gets
#    a custom date/time number format)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("icu_beds_available")
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "hospital_id"
$ws.Cells.Item(1,3).Value = "beds_in_use"
$ws.Cells.Item(1,4).Value = "vents_in_use"
$ws.Cells.Item(1,5).Value = "recorded_at"

$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = 1
$ws.Cells.Item(2,5).Value = 43932.11472222222
$ws.Cells.Item(2,5).NumberFormat = "yyyy-mm-dd h:mm:ss"

# ---------------------------------------------------------------------------
# 5. Populate "patient_responses"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("patient_responses")
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "patient_id"
$ws.Cells.Item(1,3).Value = "question_id"
$ws.Cells.Item(1,4).Value = "answer_id"
$ws.Cells.Item(1,5).Value = "recorded_at"

# ---------------------------------------------------------------------------
# 6. Populate "patients"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("patients")
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "hospital_id"
$ws.Cells.Item(1,3).Value = "patient_name"
$ws.Cells.Item(1,4).Value = "patient_mr_no"
$ws.Cells.Item(1,5).Value = "entry_point"
$ws.Cells.Item(1,6).Value = "recorded_at"

# ---------------------------------------------------------------------------
# 7. Update "questions" (was the blank default "Sheet"): rewrite its header
#    row, matching the analogous "answers" lookup table's shape.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("questions")
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "question"
$ws.Cells.Item(1,3).Value = "question_type"
$ws.Cells.Item(1,4).Value = "weight"

# ---------------------------------------------------------------------------
# 8. "answers" tab itself keeps its existing header row
#    (id, question_id, answers, weight) -- nothing to change there.
# ---------------------------------------------------------------------------
